# Adds "should " in front of the modal verb ("have"/"be on"/"allow"/"provide")
# for each requirement bullet in the networking task, per the commit
# "ADMIN: partial task and rubric for networking".

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1. "...sets of Lab Computers that have access to" -> "...that should have access to"
Replace-Text "that have" "that should have"

# 2. "BOYD (student laptops) have access to" -> "... should have access to"
Replace-Text "BOYD (student laptops) have access to" "BOYD (student laptops) should have access to"

# 3. "Teacher devices" -> "Teacher devices should have"
Replace-Text "Teacher devices" "Teacher devices should have"

# 4. "The Cyber range has no access to external networks"
#    -> "The Cyber range should  have no access to external networks" (double space preserved)
Replace-Text "The Cyber range has no access to external networks" "The Cyber range should  have no access to external networks"

# 5. "Lab, Teacher, and BOYD must be on " -> "Lab, Teacher, and BOYD should be on "
Replace-Text "Lab, Teacher, and BOYD must be on " "Lab, Teacher, and BOYD should be on "

# 6. "...However, you must allow for as much expansion..." -> "...you should allow..."
Replace-Text "However, you must allow for as much expansion as possible." "However, you should allow for as much expansion as possible."

# 7. "Provides DNS for internal domains for the cyber range" -> "Should provide DNS ..."
Replace-Text "Provides DNS for internal domains for the cyber range" "Should provide DNS for internal domains for the cyber range"

# 8. "Provides DHCP for all devices connecting t" -> "Should provide DHCP for all devices connecting t"
Replace-Text "Provides DHCP for all devices connecting t" "Should provide DHCP for all devices connecting t"
